# Continuing progress on the "inside electronics" (LEDs) section of the
# rough dimensions calculations sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "LEDs" section header (row 18) ---------------------------------
# Bold, black text - matches the style used for the other section headers
# such as "CADed Blade:" (A13) / "thick/wide/long" labels (row 16).
$ws.Range("A18").Value = "LEDs"
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").Font.Color = 0

# --- Measurement values (row 19) and labels (row 20) ---------------------
# Values entered in the same order the author originally typed them in
# (A19, B19, then the C20 label, then C19, D19 values, then D20 label),
# so new entries line up the same way in the workbook's shared string table.
$ws.Range("A19").Value = "2.13mm"
$ws.Range("B19").Value = "12mm"
$ws.Range("C20").Value = "one led piece"
$ws.Range("C19").Value = "6.9mm"
$ws.Range("D19").Value = "1.92mm"
$ws.Range("D20").Value = "between leds"

$ws.Range("A20").Value = "thick"
$ws.Range("B20").Value = "wide"

# Regular (non-bold) black text, like the other measurement-value cells
# (A15, B15, C15, F15 ...).
$valueRange19 = $ws.Range("A19:D19")
$valueRange19.Font.Bold = $false
$valueRange19.Font.Color = 0

# Last two value cells on this row are right aligned (new formatting).
$ws.Range("C19:D19").HorizontalAlignment = -4152

# Bold, black text, matching the existing "thick / wide / long" label row.
$labelRange20 = $ws.Range("A20:D20")
$labelRange20.Font.Bold = $true
$labelRange20.Font.Color = 0

# Last label cell is right aligned, matching the value column above it.
$ws.Range("D20").HorizontalAlignment = -4152

# Keep the new rows the same height as the rest of the sheet (15.75pt,
# matching the sheet's default row height).
$ws.Rows.Item(18).RowHeight = 15.75
$ws.Rows.Item(19).RowHeight = 15.75
$ws.Rows.Item(20).RowHeight = 15.75

# --- Update the active cell selection, as recorded in the saved file -----
$ws.Range("G20").Select() | Out-Null
